$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts TestCase Name/Execute Flag/
# Test case detailed description one column to the right).
$ws.Columns.Item(1).Insert()

# New column A = "Sl. No"
$ws.Range("A1").Value = "Sl. No"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Append ".java" to the TestCase Name values (now in column B)
$ws.Range("B2").Value = "Core.AlertsandNotifications.AlertsPlaceholderManagement_TestClass.java"
$ws.Range("B3").Value = "Core.AlertsandNotifications.AlertsTemplateManagement_TestClass.java"
$ws.Range("B4").Value = "Core.AlertsandNotifications.AlertsNotificationManagement_TestClass.java"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 5.54296875
$ws.Columns.Item(3).ColumnWidth = 11.7265625
$ws.Columns.Item(4).ColumnWidth = 72.7265625

$ws.Range("D8").Select()
